# Update cryptos list: refresh Price (D) and Volume(1h) (E) values,
# and fix the ARBITRUM/Aave row ordering (rows 38-39).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.859.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.74%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.841.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.55%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.620'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.96%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.82'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.15%  '

$ws.Range("E9").Value = '  +1.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0685'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0986'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.104.65'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.868.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.673'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '34.857.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0786'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.42%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.40%  '

$ws.Range("E27").Value = '  +2.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.09%  '

$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0554'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("E32").Value = '  -3.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.696'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.18%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.83%  '

$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '91.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.340.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0194'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.64'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("E43").Value = '  -1.42%  '

$ws.Range("E44").Value = '  -2.87%  '

$ws.Range("E45").Value = '  +0.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.77%  '

$ws.Range("E47").Value = '  +1.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.018.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.34%  '

$ws.Range("E49").Value = '  +5.14%  '

$ws.Range("E50").Value = '  -0.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +15.66%  '
